$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = "183096-7"
$ws.Range("B32").Value = "Clio - Greek Yogurt Bar Strawberry"
$ws.Range("C32").Value = "'2"
$ws.Range("D32").Value = "'15.45"
$ws.Range("E32").Value = "'30.90"

$ws.Range("A33").Value = "183090-0"
$ws.Range("B33").Value = "Clio - Greek Yogurt Bar Vanilla"
$ws.Range("C33").Value = "'2"
$ws.Range("D33").Value = "'15.45"
$ws.Range("E33").Value = "'30.90"
